# Applies the upstream data refresh to the "Avverkningsanmälningar" sheet:
#   - rows 2-47 are re-ordered (the source feed re-sorted the underlying
#     records between scrapes), carrying every column (A:R values plus the
#     S:Z HYPERLINK formulas) along with each record
#   - column C ("Förändrad") is bumped from 46063 to 46064 for every data
#     row, reflecting the new scrape timestamp

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row permutation: target row -> source row (both 1-based sheet rows).
# Rows that are not listed keep their own data (identity mapping), which is
# true for every row already present here since the table is 2..47.
$mapping = @{
    2=2;   3=3;   4=4;   5=6;   6=8;   7=5;   8=7;   9=9;   10=10; 11=11;
    12=12; 13=13; 14=14; 15=15; 16=16; 17=34; 18=18; 19=20; 20=21; 21=41;
    22=35; 23=46; 24=47; 25=32; 26=26; 27=31; 28=37; 29=25; 30=24; 31=27;
    32=45; 33=44; 34=29; 35=30; 36=38; 37=42; 38=19; 39=33; 40=28; 41=22;
    42=36; 43=39; 44=40; 45=23; 46=17; 47=43
}

$firstRow = 2
$lastRow  = 47
$rowCount = $lastRow - $firstRow + 1

# Snapshot the plain values (A:R) and the hyperlink formulas (S:Z)
# separately, since .Formula is needed to preserve the HYPERLINK() formulas
# instead of collapsing them to their cached text result.
$valRange  = $ws.Range("A$firstRow`:R$lastRow")
$vals      = $valRange.Value2
$formRange = $ws.Range("S$firstRow`:Z$lastRow")
$forms     = $formRange.Formula

$newVals  = New-Object 'object[,]' $rowCount,18
$newForms = New-Object 'object[,]' $rowCount,8

for ($targetRow = $firstRow; $targetRow -le $lastRow; $targetRow++) {
    $srcRow = $mapping[$targetRow]
    $ti = $targetRow - $firstRow       # 0-based row index in the new arrays
    $si = $srcRow - $firstRow          # 0-based row index in the snapshot

    for ($c = 1; $c -le 18; $c++) {
        $newVals[$ti, $c-1] = $vals[$si+1, $c]
    }
    for ($c = 1; $c -le 8; $c++) {
        $newForms[$ti, $c-1] = $forms[$si+1, $c]
    }
}

$valRange.Value2  = $newVals
$formRange.Formula = $newForms

# Bump the "Förändrad" column to the new scrape date for every data row.
$ws.Range("C$firstRow`:C$lastRow").Value2 = 46064

# Re-writing a wrapped multi-line cell (row 2's species list, the only cell
# with an embedded line break) makes Excel auto-fit that row taller; restore
# the original fixed row height so the sheet layout matches the source
# feed's formatting exactly. Other rows are left alone since they were never
# auto-fit and already keep their original height.
$ws.Rows.Item(2).RowHeight = 15
